# Natmi following Dr Hou advice
# Add the "ECs" sending-cluster block (Efna5 -> Epha1 onto ECs/FAPs/M2/sCs)
# ahead of the existing FAPs/sCs blocks, and refresh all of the
# expression/specificity statistics for every Sending/Target cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row: A..T (20 columns).
$rows = @(
    @("ECs",  "Efna5", "Epha1", "ECs",  1, 0.3333333333333333, 0.09353433333333333, 0.280603,   0.0471360296668244, 0.0471360296668244, 3, 1, 3.746298666666667, 11.238896,  0.2944715732214049, 0.294471573221405,  0.3504075482542222, 3.153667934288,    0.0138802208114006,  0.0138802208114006),
    @("ECs",  "Efna5", "Epha1", "FAPs", 1, 0.3333333333333333, 0.09353433333333333, 0.280603,   0.0471360296668244, 0.0471360296668244, 3, 1, 4.370261666666667, 13.110785,  0.3435171466234404, 0.3435171466234404, 0.4087695114838889, 3.678925603355,    0.01619203441430535, 0.01619203441430535),
    @("ECs",  "Efna5", "Epha1", "M2",   1, 0.3333333333333333, 0.09353433333333333, 0.280603,   0.0471360296668244, 0.0471360296668244, 3, 1, 1.837969333333333, 5.513908,   0.1444705212467569, 0.1444705212467569, 0.1719132362804444, 1.547219126524,    0.006809766775468716, 0.006809766775468716),
    @("ECs",  "Efna5", "Epha1", "sCs",  1, 0.3333333333333333, 0.09353433333333333, 0.280603,   0.0471360296668244, 0.0471360296668244, 3, 1, 2.767576666666667, 8.30273,    0.2175407589083977, 0.2175407589083977, 0.2588634384655555, 2.32977094619,     0.01025400766564973, 0.01025400766564973),
    @("FAPs", "Efna5", "Epha1", "ECs",  3, 1,                  1.367717666666667,   4.103153,   0.6892525793933763, 0.6892525793933763, 3, 1, 3.746298666666667, 11.238896,  0.2944715732214049, 0.294471573221405,  5.123878871009778,  46.114909839088,   0.2029652914008788,  0.2029652914008789),
    @("FAPs", "Efna5", "Epha1", "FAPs", 3, 1,                  1.367717666666667,   4.103153,   0.6892525793933763, 0.6892525793933763, 3, 1, 4.370261666666667, 13.110785,  0.3435171466234404, 0.3435171466234404, 5.977284089456111,  53.795556805105,   0.2367700793760589,  0.2367700793760589),
    @("FAPs", "Efna5", "Epha1", "M2",   3, 1,                  1.367717666666667,   4.103153,   0.6892525793933763, 0.6892525793933763, 3, 1, 1.837969333333333, 5.513908,   0.1444705212467569, 0.1444705212467569, 2.513823127991555,  22.624408151924,   0.09957667941563272, 0.09957667941563272),
    @("FAPs", "Efna5", "Epha1", "sCs",  3, 1,                  1.367717666666667,   4.103153,   0.6892525793933763, 0.6892525793933763, 3, 1, 2.767576666666667, 8.30273,    0.2175407589083977, 0.2175407589083977, 3.785263500854445,  34.06737150769,    0.1499405292008057,  0.1499405292008057),
    @("sCs",  "Efna5", "Epha1", "ECs",  3, 1,                  0.523097,             1.569291,   0.2636113909397994, 0.2636113909397994, 3, 1, 3.746298666666667, 11.238896,  0.2944715732214049, 0.294471573221405,  1.959677593637333,  17.637098342736,   0.07762606100912554, 0.07762606100912556),
    @("sCs",  "Efna5", "Epha1", "FAPs", 3, 1,                  0.523097,             1.569291,   0.2636113909397994, 0.2636113909397994, 3, 1, 4.370261666666667, 13.110785,  0.3435171466234404, 0.3435171466234404, 2.286070767048334,  20.574636903435,   0.09055503283307614, 0.09055503283307614),
    @("sCs",  "Efna5", "Epha1", "M2",   3, 1,                  0.523097,             1.569291,   0.2636113909397994, 0.2636113909397994, 3, 1, 1.837969333333333, 5.513908,   0.1444705212467569, 0.1444705212467569, 0.9614362443586667, 8.652926199228,    0.03808407505565542, 0.03808407505565542),
    @("sCs",  "Efna5", "Epha1", "sCs",  3, 1,                  0.523097,             1.569291,   0.2636113909397994, 0.2636113909397994, 3, 1, 2.767576666666667, 8.30273,    0.2175407589083977, 0.2175407589083977, 1.447711051603334,  13.02939946443,    0.05734622204194229, 0.05734622204194229)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T13").Value = $data
